$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 308.33334
$ws.Range("I33").Value = 280.3125
$ws.Range("K33").Value = 280.3125
$ws.Range("M33").Value = -51.3125
# Row 92
$ws.Range("H92").Value = 1254.1666
$ws.Range("I92").Value = 1181.48
$ws.Range("K92").Value = 1181.48
$ws.Range("M92").Value = 66.51999999999998
# Row 100
$ws.Range("H100").Value = 4737.5884
$ws.Range("I100").Value = 3034.5715
$ws.Range("J100").Value = 5929.7
$ws.Range("K100").Value = 3034.5715
$ws.Range("L100").Value = 5929.7
$ws.Range("M100").Value = -2493.5715
$ws.Range("N100").Value = -7011.7
# Row 132
$ws.Range("H132").Value = 1870.2
$ws.Range("I132").Value = 1629.7838
$ws.Range("K132").Value = 4889.3514
$ws.Range("M132").Value = -2359.3514
# Row 133
$ws.Range("H133").Value = 103844.625
$ws.Range("J133").Value = 103844.625
$ws.Range("L133").Value = 103844.625
$ws.Range("N133").Value = -113964.625
# Row 135
$ws.Range("H135").Value = 934
$ws.Range("I135").Value = 934
$ws.Range("K135").Value = 8406
$ws.Range("M135").Value = -5871

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4038.4583
$ws.Range("I2").Value = 3819.1765
$ws.Range("J2").Value = 4571
$ws.Range("K2").Value = 3819.1765
$ws.Range("L2").Value = 4571
$ws.Range("M2").Value = -3706.1765
$ws.Range("N2").Value = -4797
# Row 32
$ws.Range("H32").Value = 95234.35000000001
$ws.Range("I32").Value = 121592.51
$ws.Range("J32").Value = 14277.143
$ws.Range("K32").Value = 121592.51
$ws.Range("L32").Value = 14277.143
$ws.Range("M32").Value = -121305.51
$ws.Range("N32").Value = -14851.143
# Row 61
$ws.Range("H61").Value = 5339.3335
$ws.Range("I61").Value = 4148.08
$ws.Range("K61").Value = 4148.08
$ws.Range("M61").Value = -3936.08
# Row 63
$ws.Range("H63").Value = 6410.5713
$ws.Range("I63").Value = 2535.5715
$ws.Range("J63").Value = 10285.571
$ws.Range("K63").Value = 2535.5715
$ws.Range("L63").Value = 10285.571
$ws.Range("M63").Value = -1849.5715
$ws.Range("N63").Value = -11657.571
# Row 66
$ws.Range("H66").Value = 6410.5713
$ws.Range("I66").Value = 2535.5715
$ws.Range("J66").Value = 10285.571
$ws.Range("K66").Value = 12677.8575
$ws.Range("L66").Value = 51427.855
$ws.Range("M66").Value = -9245.8575
$ws.Range("N66").Value = -58291.855
# Row 116
$ws.Range("H116").Value = 4038.4583
$ws.Range("I116").Value = 3819.1765
$ws.Range("J116").Value = 4571
$ws.Range("K116").Value = 3819.1765
$ws.Range("L116").Value = 4571
$ws.Range("M116").Value = -1525.1765
$ws.Range("N116").Value = -9159
# Row 130
$ws.Range("H130").Value = 66741.5
$ws.Range("J130").Value = 66741.5
$ws.Range("L130").Value = 66741.5
$ws.Range("N130").Value = -76781.5
# Row 136
$ws.Range("H136").Value = 5339.3335
$ws.Range("I136").Value = 4148.08
$ws.Range("K136").Value = 12444.24
$ws.Range("M136").Value = -9894.24

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4038.4583
$ws.Range("I3").Value = 3819.1765
$ws.Range("J3").Value = 4571
$ws.Range("K3").Value = 3819.1765
$ws.Range("L3").Value = 4571
$ws.Range("M3").Value = -3705.1765
$ws.Range("N3").Value = -4799
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# Row 54
$ws.Range("H54").Value = 6367.643
$ws.Range("I54").Value = 2242.6924
$ws.Range("K54").Value = 2242.6924
$ws.Range("M54").Value = -1758.6924
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 99
$ws.Range("H99").Value = 4739.909
$ws.Range("I99").Value = 4573.1665
$ws.Range("J99").Value = 4940
$ws.Range("K99").Value = 4573.1665
$ws.Range("L99").Value = 4940
$ws.Range("M99").Value = -3075.1665
$ws.Range("N99").Value = -7936
# Row 105
$ws.Range("H105").Value = 10062.9
$ws.Range("I105").Value = 13905.454
$ws.Range("J105").Value = 5366.4443
$ws.Range("K105").Value = 13905.454
$ws.Range("L105").Value = 5366.4443
$ws.Range("M105").Value = -12158.454
$ws.Range("N105").Value = -8860.444299999999
# Row 107
$ws.Range("H107").Value = 3302.0425
$ws.Range("I107").Value = 2277.6943
$ws.Range("J107").Value = 6654.4546
$ws.Range("K107").Value = 2277.6943
$ws.Range("L107").Value = 6654.4546
$ws.Range("M107").Value = -357.6943000000001
$ws.Range("N107").Value = -10494.4546
# Row 134
$ws.Range("H134").Value = 3257.111
$ws.Range("I134").Value = 3374.923
$ws.Range("K134").Value = 10124.769
$ws.Range("M134").Value = -7589.769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4400.1665
$ws.Range("I31").Value = 2140.182
$ws.Range("J31").Value = 6312.4614
$ws.Range("K31").Value = 2140.182
$ws.Range("L31").Value = 6312.4614
$ws.Range("M31").Value = -1845.182
$ws.Range("N31").Value = -6902.4614
# Row 34
$ws.Range("H34").Value = 4400.1665
$ws.Range("I34").Value = 2140.182
$ws.Range("J34").Value = 6312.4614
$ws.Range("K34").Value = 2140.182
$ws.Range("L34").Value = 6312.4614
$ws.Range("M34").Value = -1938.182
$ws.Range("N34").Value = -6716.4614
# Row 132
$ws.Range("H132").Value = 1324.12
$ws.Range("I132").Value = 982.7826
$ws.Range("K132").Value = 2948.3478
$ws.Range("M132").Value = -418.3478
# Row 134
$ws.Range("H134").Value = 23116.2
$ws.Range("I134").Value = 24453.072
$ws.Range("J134").Value = 4400
$ws.Range("K134").Value = 73359.216
$ws.Range("L134").Value = 13200
$ws.Range("M134").Value = -70824.216
$ws.Range("N134").Value = -18270

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 1285.5714
$ws.Range("J80").Value = 1524.75
$ws.Range("L80").Value = 4574.25
$ws.Range("N80").Value = -6446.25
# Row 83
$ws.Range("H83").Value = 1285.5714
$ws.Range("J83").Value = 1524.75
$ws.Range("L83").Value = 13722.75
$ws.Range("N83").Value = -23082.75
# Row 92
$ws.Range("H92").Value = 135
$ws.Range("I92").Value = 131.25
$ws.Range("J92").Value = 138.75
$ws.Range("K92").Value = 393.75
$ws.Range("L92").Value = 416.25
$ws.Range("M92").Value = 854.25
$ws.Range("N92").Value = -2912.25
# Row 131
$ws.Range("H131").Value = 3180525
$ws.Range("I131").Value = 1850
$ws.Range("J131").Value = 3515122.2
$ws.Range("K131").Value = 5550
$ws.Range("L131").Value = 10545366.6
$ws.Range("M131").Value = -510
$ws.Range("N131").Value = -10555446.6
# Row 139
$ws.Range("H139").Value = 3043.5
$ws.Range("I139").Value = 3043.5
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 9130.5
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -3990.5
$ws.Range("N139").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 33338136
$ws.Range("I70").Value = 333333340
$ws.Range("J70").Value = 5334
$ws.Range("K70").Value = 333333340
$ws.Range("L70").Value = 5334
$ws.Range("M70").Value = -333333070
$ws.Range("N70").Value = -5874
# Row 73
$ws.Range("H73").Value = 33338136
$ws.Range("I73").Value = 333333340
$ws.Range("J73").Value = 5334
$ws.Range("K73").Value = 333333340
$ws.Range("L73").Value = 5334
$ws.Range("M73").Value = -333332404
$ws.Range("N73").Value = -7206
# Row 107
$ws.Range("H107").Value = 2375.6667
$ws.Range("I107").Value = 1090
$ws.Range("J107").Value = 2632.8
$ws.Range("K107").Value = 1090
$ws.Range("L107").Value = 2632.8
$ws.Range("M107").Value = 830
$ws.Range("N107").Value = -6472.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 123855.95
$ws.Range("I2").Value = 139552.88
$ws.Range("J2").Value = 1420
$ws.Range("K2").Value = 139552.88
$ws.Range("L2").Value = 1420
$ws.Range("M2").Value = -139440.88
$ws.Range("N2").Value = -1644
# Row 46
$ws.Range("H46").Value = 1398.7941
$ws.Range("I46").Value = 921.3125
$ws.Range("J46").Value = 1823.2222
$ws.Range("K46").Value = 921.3125
$ws.Range("L46").Value = 1823.2222
$ws.Range("M46").Value = -733.3125
$ws.Range("N46").Value = -2199.2222
# Row 55
$ws.Range("H55").Value = 2177.0908
$ws.Range("I55").Value = 1564.5714
$ws.Range("K55").Value = 1564.5714
$ws.Range("M55").Value = -1391.5714

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 11421.105
$ws.Range("J62").Value = 12846.23
$ws.Range("L62").Value = 12846.23
$ws.Range("N62").Value = -14094.23
# Row 65
$ws.Range("H65").Value = 11421.105
$ws.Range("J65").Value = 12846.23
$ws.Range("L65").Value = 64231.14999999999
$ws.Range("N65").Value = -70471.14999999999
# Row 100
$ws.Range("H100").Value = 948.75
$ws.Range("I100").Value = 398.33334
$ws.Range("J100").Value = 2600
$ws.Range("K100").Value = 796.66668
$ws.Range("L100").Value = 5200
$ws.Range("M100").Value = -255.66668
$ws.Range("N100").Value = -6282
# Row 113
$ws.Range("H113").Value = 938.65515
$ws.Range("I113").Value = 798.7619
$ws.Range("J113").Value = 1305.875
$ws.Range("K113").Value = 2396.2857
$ws.Range("L113").Value = 3917.625
$ws.Range("M113").Value = -226.2856999999999
$ws.Range("N113").Value = -8257.625
# Row 136
$ws.Range("H136").Value = 4056.4736
$ws.Range("I136").Value = 1175.125
$ws.Range("J136").Value = 6152
$ws.Range("K136").Value = 3525.375
$ws.Range("L136").Value = 18456
$ws.Range("M136").Value = -975.375
$ws.Range("N136").Value = -23556
